# Updates cryptos list price (D) and volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.714.24"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").Value = "'1.638.08"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'217.59"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("E6").Value = "  -0.99%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("E9").Value = "  -0.81%  "

$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("E11").Value = "  +0.09%  "

$ws.Range("D12").Value = "'1.865.11"
$ws.Range("E12").Value = "  -0.84%  "

$ws.Range("D13").Value = "'1.642.06"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").Value = "'0.525"
$ws.Range("E15").Value = "  -1.51%  "

$ws.Range("D16").Value = "'64.45"
$ws.Range("E16").Value = "  -1.70%  "

$ws.Range("D17").Value = "'26.677.69"
$ws.Range("E17").Value = "  -0.59%  "

$ws.Range("E18").Value = "  -2.48%  "

$ws.Range("D19").Value = "'211.33"
$ws.Range("E19").Value = "  -3.47%  "

$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("E21").Value = "  -0.81%  "

$ws.Range("D22").Value = "'6.19"
$ws.Range("E22").Value = "  -1.49%  "

$ws.Range("D23").Value = "'2.31"
$ws.Range("E23").Value = "  -3.09%  "

$ws.Range("E24").Value = "  -2.82%  "

$ws.Range("D25").Value = "'146.60"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("E27").Value = "  -1.99%  "

$ws.Range("D28").Value = "'7.07"
$ws.Range("E28").Value = "  -0.76%  "

$ws.Range("D29").Value = "'15.55"
$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("D30").Value = "'0.0502"
$ws.Range("E30").Value = "  -2.91%  "

$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("D32").Value = "'3.35"
$ws.Range("E32").Value = "  -0.28%  "

$ws.Range("E33").Value = "  -1.23%  "

$ws.Range("D34").Value = "'1.268.17"
$ws.Range("E34").Value = "  -1.44%  "

$ws.Range("D35").Value = "'1.53"
$ws.Range("E35").Value = "  -1.22%  "

$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("E37").Value = "  -2.34%  "

$ws.Range("E38").Value = "  -2.00%  "

$ws.Range("E39").Value = "  -2.98%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("E41").Value = "  -1.52%  "

$ws.Range("E42").Value = "  -2.84%  "

$ws.Range("E43").Value = "  -3.76%  "

$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("D45").Value = "'91.48"
$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("D46").Value = "'60.17"
$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("E47").Value = "  -1.76%  "

$ws.Range("E48").Value = "  +0.33%  "

$ws.Range("E49").Value = "  -3.05%  "

$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("D51").Value = "'0.407"
$ws.Range("E51").Value = "  -0.37%  "
